# Update argument-description text (ConvertTypes, Delimiter, HeaderRowNum,
# ShowMissingsAs, Encoding, DecimalSeparator, HeaderRow, QuoteAllStrings, EOL)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = "Controls whether fields in the file are converted to typed values or remain as strings, and sets the treatment of `"quoted fields`" and space characters.`r`n`r`nConvertTypes should be a string of zero or more letters from allowed characters ``NDBETQ``.`r`n`r`nThe most commonly useful letters are:`r`n1) ``N`` number fields are returned as numbers (Doubles).`r`n2) ``D`` date fields (that respect DateFormat) are returned as Dates.`r`n3) ``B`` fields matching TrueStrings or FalseStrings are returned as Booleans.`r`n`r`nConvertTypes is optional and defaults to the null string for no type conversion. ``TRUE`` is equivalent to ``NDB`` and ``FALSE`` to the null string.`r`n`r`nThree further options are available:`r`n4) ``E`` fields that match Excel errors are converted to error values. There are fourteen of these, including ``#N/A``, ``#NAME?``, ``#VALUE!`` and ``#DIV/0!``.`r`n5) ``T`` leading and trailing spaces are trimmed from fields. In the case of quoted fields, this will not remove spaces between the quotes.`r`n6) ``Q`` conversion happens for both quoted and unquoted fields; otherwise only unquoted fields are converted.`r`n`r`nFor most files, correct type conversion can be achieved with ConvertTypes as a string which applies for all columns, but type conversion can also be specified on a per-column basis.`r`n`r`nEnter an array (or range) with two columns or two rows, column numbers on the left/top and type conversion (subset of ``NDBETQ``) on the right/bottom. Instead of column numbers, you can enter strings matching the contents of the header row, and a column number of zero applies to all columns not otherwise referenced.`r`n`r`nFor convenience when calling from VBA, you can pass an array of two element arrays such as ``Array(Array(0,`"N`"),Array(3,`"`"),Array(`"Phone`",`"`"))`` to convert all numbers in a file into numbers in the return except for those in column 3 and in the column(s) headed `"Phone`"."
$ws.Range("I2").Value = "By default, CSVRead will try to detect a file's delimiter as the first instance of comma, tab, semi-colon, colon or pipe found outside quoted regions in the first 10,000 characters of the file. If it can't auto-detect the delimiter, it will assume comma. If your file includes a different character or string delimiter you should pass that as the Delimiter argument.`r`n`r`nAlternatively, enter ``FALSE`` as the delimiter to treat the file as `"not a delimited file`". In this case the return will mimic how the file would appear in a text editor such as NotePad. The file will be split into lines at all line breaks (irrespective of double quotes) and each element of the return will be a line of the file."
$ws.Range("S2").Value = "The row in the file containing headers. Type conversion is not applied to fields in the header row, though leading and trailing spaces are trimmed.`r`n`r`nThis argument is most useful when calling from VBA, with SkipToRow set to one more than HeaderRowNum. In that case the function returns the rows starting from SkipToRow, and the header row is returned via the by-reference argument HeaderRow. Optional and defaults to 0."
$ws.Range("AI2").Value = "Fields which are missing in the file (consecutive delimiters) or match one of the MissingStrings are returned in the array as ShowMissingsAs. Defaults to Empty, but the null string or ``#N/A!`` error value can be good alternatives.`r`n`r`nIf NumRows is greater than the number of rows in the file then the return is `"padded`" with the value of ShowMissingsAs. Likewise, if NumCols is greater than the number of columns in the file."
$ws.Range("AK2").Value = "Allowed entries are ``ASCII``, ``ANSI``, ``UTF-8``, or ``UTF-16``. For most files this argument can be omitted and CSVRead will detect the file's encoding. If auto-detection does not work, then it's possible that the file is encoded ``UTF-8`` or ``UTF-16`` but without a byte option mark to identify the encoding. Experiment with Encoding as each of ``UTF-8`` and ``UTF-16``."
$ws.Range("AM2").Value = "In many places in the world, floating point number decimals are separated with a comma instead of a period (3,14 vs. 3.14). CSVRead can correctly parse these numbers by passing in the DecimalSeparator as a comma, in which case comma ceases to be a candidate if the parser needs to guess the Delimiter."
$ws.Range("AO2").Value = "This by-reference argument is for use from VBA (as opposed to from Excel formulas). It is populated with the contents of the header row, with no type conversion, though leading and trailing spaces are removed."
$ws.Range("I3").Value = "If ``TRUE`` (the default) then elements of Data that are strings are quoted before being written to file, other elements (Numbers, Booleans, Errors) are not quoted. If ``FALSE`` then the only elements of Data that are quoted are strings containing Delimiter, line feed, carriage return or double quote. In all cases, double quotes are escaped by another double quote."
$ws.Range("S3").Value = "Controls the line endings of the file written. Enter ``Windows`` (the default), ``Unix`` or ``Mac``. Also supports the line-ending characters themselves (ascii 13 + ascii 10, ascii 10, ascii 13) or the strings ``CRLF``, ``LF`` or ``CR``. The last line of the file is written with a line ending."

# Column widths follow Excel's own best-fit recalculation once the
# descriptions above changed length; set them to match.
$ws.Columns.Item(7).ColumnWidth = 109.9169
$ws.Columns.Item(35).ColumnWidth = 39.7527
$ws.Columns.Item(37).ColumnWidth = 38.4184
$ws.Columns.Item(39).ColumnWidth = 37.5863
$ws.Columns.Item(41).ColumnWidth = 37.7527

# Row heights grow because of the longer, wrapped text above.
$ws.Rows.Item(2).RowHeight = 409.5
$ws.Rows.Item(3).RowHeight = 114

